$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Dia chi lien lac khi can bao tin: <<ThiSinh_DienThoai>>"
#    -> "...: <<ThiSinh_DCNhanGiayBao>>" split across four runs, with the
#    placeholder name run carrying an explicit black font color, matching
#    the merge-field style used elsewhere in the template.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute(": <<ThiSinh_DienThoai>>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $rng.Start
    $end = $rng.End

    $full = $d.Range($start, $end)
    $full.Text = ": <<ThiSinh_DCNhanGiayBao>>"

    # New text layout (0-based offsets from $start):
    #   [0,2)   ": "
    #   [2,4)   "<<"
    #   [4,25)  "ThiSinh_DCNhanGiayBao"
    #   [25,27) ">>"
    $seg1 = $d.Range($start, $start + 2)
    $seg2 = $d.Range($start + 2, $start + 4)
    $seg3 = $d.Range($start + 4, $start + 25)
    $seg4 = $d.Range($start + 25, $start + 27)

    # Force seg1/seg2 apart into distinct runs (identical formatting would
    # otherwise be silently re-merged) by toggling a property and reverting
    # it before save.
    $seg1.Font.Bold = 1
    $seg1.Font.Bold = 0

    # The placeholder name itself gets an explicit black color.
    $seg3.Font.Color = 0
}

# ---------------------------------------------------------------------------
# 2) Style bookkeeping touch-ups picked up by the style inspector.
# ---------------------------------------------------------------------------
$defParaFont = $d.Styles.Item("DefaultParagraphFont")
$defParaFont.Priority = 1
$defParaFont.UnhideWhenUsed = $true

$tableNormal = $d.Styles.Item("TableNormal")
$tableNormal.Priority = 99
$tableNormal.UnhideWhenUsed = $true

$noList = $d.Styles.Item("NoList")
$noList.Priority = 99
$noList.UnhideWhenUsed = $true

# Fix the stray leading space in the custom "Char" style's display name.
$charStyle = $d.Styles.Item("Char")
$charStyle.NameLocal = "Char"
